$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9: only the taxon sort order (B) changes.
$ws.Range("B9").Value = 96720

# Row 10 <-> Row 12: the species-identifying columns (A,B,D,E,F,G,H,Q,R) get
# swapped between these two rows. Row 12 additionally carries an empty
# placeholder cell in AF that moves to row 10 along with the rest of the
# species data.
$ws.Range("AF12").Cut($ws.Range("AF10")) | Out-Null

$ws.Range("A10").Value = 112044164
$ws.Range("B10").Value = 89058
$ws.Range("D10").Value = "LC"
$ws.Range("E10").Value = 256703
$ws.Range("F10").Value = "Tallfingersvamp"
$ws.Range("G10").Value = "Ramaria eosanguinea"
$ws.Range("H10").Value = "R.H.Petersen"
$ws.Range("Q10").Value = 554725
$ws.Range("R10").Value = 6697591

$ws.Range("A12").Value = 112044172
$ws.Range("B12").Value = 90792
$ws.Range("D12").Value = "NT"
$ws.Range("E12").Value = 4361
$ws.Range("F12").Value = "Orange taggsvamp"
$ws.Range("G12").Value = "Hydnellum aurantiacum"
$ws.Range("H12").Value = "(Batsch:Fr.) P.Karst."
$ws.Range("Q12").Value = 554722
$ws.Range("R12").Value = 6697604

# Row 11: only the taxon sort order (B) changes.
$ws.Range("B11").Value = 96720

# Rows 23-28: species-identifying columns get reshuffled between the rows
# (A,B,D,E,F,G,H,Q,R for 23/25/27/28; A,B,Q,R only for 24/26 since their
# species data is identical before and after).
$ws.Range("A23").Value = 112044179
$ws.Range("B23").Value = 96720
$ws.Range("D23").Value = "VU"
$ws.Range("E23").Value = 220787
$ws.Range("F23").Value = "Knärot"
$ws.Range("G23").Value = "Goodyera repens"
$ws.Range("H23").Value = "(L.) R. Br."
$ws.Range("Q23").Value = 554795
$ws.Range("R23").Value = 6697596

$ws.Range("A24").Value = 112044168
$ws.Range("B24").Value = 89979
$ws.Range("Q24").Value = 554761
$ws.Range("R24").Value = 6697614

$ws.Range("A25").Value = 112044155
$ws.Range("B25").Value = 89539
$ws.Range("D25").Value = "NT"
$ws.Range("E25").Value = 1202
$ws.Range("F25").Value = "Ullticka"
$ws.Range("G25").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H25").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("Q25").Value = 554761
$ws.Range("R25").Value = 6697629

$ws.Range("A26").Value = 112044169
$ws.Range("B26").Value = 89979
$ws.Range("Q26").Value = 554765
$ws.Range("R26").Value = 6697617

$ws.Range("A27").Value = 112044157
$ws.Range("B27").Value = 89539
$ws.Range("D27").Value = "NT"
$ws.Range("E27").Value = 1202
$ws.Range("F27").Value = "Ullticka"
$ws.Range("G27").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H27").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("Q27").Value = 554764
$ws.Range("R27").Value = 6697616

$ws.Range("A28").Value = 112044180
$ws.Range("B28").Value = 96720
$ws.Range("D28").Value = "VU"
$ws.Range("E28").Value = 220787
$ws.Range("F28").Value = "Knärot"
$ws.Range("G28").Value = "Goodyera repens"
$ws.Range("H28").Value = "(L.) R. Br."
$ws.Range("Q28").Value = 554839
$ws.Range("R28").Value = 6697581

# Rows 29-30: only the taxon sort order (B) changes.
$ws.Range("B29").Value = 89539
$ws.Range("B30").Value = 89539
